$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (per diff: 14.08988764044944 -> 16.28988764044944).
# ColumnWidth is stored/quantized in whole-pixel steps (1/7 of a "character"
# unit here), so feed it the pre-image that lands on the closest reachable
# pixel-quantized width to the target.
$ws.Columns.Item(1).ColumnWidth = 15.571428571428571

# Row 2: Account Number becomes a numeric value, Netpay updated
$ws.Range("A2").Value = 32145698741
$ws.Range("C2").Value = 40989.1

# New Row 3 with a new account entry
$ws.Range("A3").Value = 123654789963
$ws.Range("B3").Value = "Vidya Sagar pogiri"
$ws.Range("C3").Value = 9793.33
$ws.Range("D3").Value = "November"
